$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (pushes old rows 13-21 down to 14-22) ---
$ws.Rows.Item(13).Insert()

# Give the newly inserted row's B/C cells the same formatting (wrap text styles)
# used by the other content cells in columns B/C, by copying formats from row 10.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)

# --- Row 10 (Objetivos:) - replace the (misplaced) teacher name with the real objectives text ---
$ws.Range("B10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão de negócios."
$ws.Range("C10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre gestão de negócios."

# --- Row 13 (new, blank row under "Docentes responsáveis:") - holds the teacher name ---
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

# --- Row 14 (Programa resumido:) - replace "Semestral" with the real short syllabus text ---
$ws.Range("B14").Value = "A definir, de acordo com o tópico programado."
$ws.Range("C14").Value = "A definir, de acordo com o tópico programado."

# --- Row 16 (Programa:) - replace leftover date with the real syllabus text ---
$ws.Range("B16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares referentes a gestão de negócios relevantes para a formação de um profissional de Engenharia."
$ws.Range("C16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares referentes a gestão de negócios relevantes para a formação de um profissional de Engenharia."

# --- Row 19 (Método:) - replace leftover teacher name with the method text ---
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# --- Row 20 (Critério:) - replace leftover method text with "Provas e trabalhos." ---
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# --- Row 21 (Norma de recuperação:) - replace leftover "Provas e trabalhos." with recovery text ---
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# --- Row 22 (Bibliografia:) - replace leftover recovery text with bibliography ---
$ws.Range("B22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."
$ws.Range("C22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas na área de gestão e produção."

# --- Split the merged column A/B width definition into independent column entries ---
# (column A keeps its original 30.71-char width, column B keeps its own 60.71-char width)
$ws.Columns.Item(1).Hidden = $false

Write-Output "done"
